$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 36885.332
$ws.Range("J3").Value = 36885.332
$ws.Range("L3").Value = 36885.332
$ws.Range("N3").Value = -37113.332
$ws.Range("H40").Value = 3483.3333
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 3483.3333
$ws.Range("K40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("M40").Value = 3483.3333
$ws.Range("N40").Value = -3833.3333
$ws.Range("H41").Value = 516.6129
$ws.Range("I41").Value = 616.5714
$ws.Range("J41").Value = 434.29413
$ws.Range("K41").Value = 616.5714
$ws.Range("L41").Value = 434.29413
$ws.Range("M41").Value = -176.5714
$ws.Range("N41").Value = -1314.29413
$ws.Range("H53").Value = 205.81818
$ws.Range("I53").Value = 266
$ws.Range("K53").Value = 266
$ws.Range("M53").Value = 371
$ws.Range("H102").Value = 36885.332
$ws.Range("J102").Value = 36885.332
$ws.Range("L102").Value = 36885.332
$ws.Range("N102").Value = -43375.332
$ws.Range("H104").Value = 648.25
$ws.Range("I104").Value = 799
$ws.Range("J104").Value = 196
$ws.Range("K104").Value = 2397
$ws.Range("L104").Value = 588
$ws.Range("M104").Value = -650
$ws.Range("N104").Value = -4082
$ws.Range("H113").Value = 4573.5
$ws.Range("I113").Value = 3433.75
$ws.Range("K113").Value = 3433.75
$ws.Range("M113").Value = -179.75
$ws.Range("H132").Value = 1171355
$ws.Range("I132").Value = 1482551.8
$ws.Range("J132").Value = 4367
$ws.Range("K132").Value = 4447655.4
$ws.Range("L132").Value = 13101
$ws.Range("M132").Value = -4445125.4
$ws.Range("N132").Value = -18161
$ws.Range("H135").Value = 1304.3334
$ws.Range("I135").Value = 1114.8334
$ws.Range("K135").Value = 10033.5006
$ws.Range("M135").Value = -7498.500599999999
$ws.Range("H138").Value = 5177.7646
$ws.Range("I138").Value = 4900
$ws.Range("J138").Value = 5214.8
$ws.Range("K138").Value = 14700
$ws.Range("L138").Value = 15644.4
$ws.Range("M138").Value = -9560
$ws.Range("N138").Value = -25924.4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3161.25
$ws.Range("I2").Value = 3161.25
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 3161.25
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -3048.25
$ws.Range("H15").Value = 2000
$ws.Range("J15").Value = 2000
$ws.Range("L15").Value = 2000
$ws.Range("N15").Value = -2700
$ws.Range("H32").Value = 126439.836
$ws.Range("I32").Value = 123708.34
$ws.Range("K32").Value = 123708.34
$ws.Range("M32").Value = -123421.34
$ws.Range("H61").Value = 3203.6875
$ws.Range("I61").Value = 3261.4285
$ws.Range("K61").Value = 3261.4285
$ws.Range("M61").Value = -3049.4285
$ws.Range("H75").Value = 23173
$ws.Range("J75").Value = 23173
$ws.Range("L75").Value = 23173
$ws.Range("N75").Value = -24921
$ws.Range("H78").Value = 23173
$ws.Range("J78").Value = 23173
$ws.Range("L78").Value = 69519
$ws.Range("N78").Value = -78255
$ws.Range("H116").Value = 3161.25
$ws.Range("I116").Value = 3161.25
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 3161.25
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -867.25
$ws.Range("H132").Value = 2741.0417
$ws.Range("I132").Value = 2594.3
$ws.Range("K132").Value = 7782.900000000001
$ws.Range("M132").Value = -5252.900000000001
$ws.Range("H136").Value = 3203.6875
$ws.Range("I136").Value = 3261.4285
$ws.Range("K136").Value = 9784.2855
$ws.Range("M136").Value = -7234.2855
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3161.25
$ws.Range("I3").Value = 3161.25
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 3161.25
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -3047.25
$ws.Range("H20").Value = 6153.0605
$ws.Range("I20").Value = 4601.6
$ws.Range("K20").Value = 4601.6
$ws.Range("M20").Value = -4354.6
$ws.Range("H99").Value = 3031.3
$ws.Range("I99").Value = 3031.3
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3031.3
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -1533.3
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2362.7778
$ws.Range("I16").Value = 2362.7778
$ws.Range("K16").Value = 2362.7778
$ws.Range("M16").Value = -2075.7778
$ws.Range("H31").Value = 2698.75
$ws.Range("I31").Value = 2698.75
$ws.Range("K31").Value = 2698.75
$ws.Range("M31").Value = -2403.75
$ws.Range("H34").Value = 2698.75
$ws.Range("I34").Value = 2698.75
$ws.Range("K34").Value = 2698.75
$ws.Range("M34").Value = -2496.75
$ws.Range("H58").Value = 9715.833000000001
$ws.Range("I58").Value = 13851.375
$ws.Range("J58").Value = 1444.75
$ws.Range("K58").Value = 13851.375
$ws.Range("L58").Value = 1444.75
$ws.Range("M58").Value = -13648.375
$ws.Range("N58").Value = -1850.75
$ws.Range("H113").Value = 2362.7778
$ws.Range("I113").Value = 2362.7778
$ws.Range("K113").Value = 2362.7778
$ws.Range("M113").Value = -192.7777999999998
$ws.Range("H122").Value = 2011.76
$ws.Range("I122").Value = 1326.0526
$ws.Range("J122").Value = 4183.1665
$ws.Range("K122").Value = 3978.1578
$ws.Range("L122").Value = 12549.4995
$ws.Range("M122").Value = -1528.1578
$ws.Range("N122").Value = -17449.4995
$ws.Range("H132").Value = 5469
$ws.Range("I132").Value = 5469
$ws.Range("K132").Value = 16407
$ws.Range("M132").Value = -13877
$ws.Range("H134").Value = 2529.3845
$ws.Range("I134").Value = 1982
$ws.Range("K134").Value = 5946
$ws.Range("M134").Value = -3411
$ws.Range("H136").Value = 9715.833000000001
$ws.Range("I136").Value = 13851.375
$ws.Range("J136").Value = 1444.75
$ws.Range("K136").Value = 41554.125
$ws.Range("L136").Value = 4334.25
$ws.Range("M136").Value = -39004.125
$ws.Range("N136").Value = -9434.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 102.5
$ws.Range("I2").Value = 86.111115
$ws.Range("K2").Value = 516.66669
$ws.Range("M2").Value = -403.66669
$ws.Range("H17").Value = 177.66667
$ws.Range("I17").Value = 127.333336
$ws.Range("J17").Value = 278.33334
$ws.Range("K17").Value = 382.000008
$ws.Range("L17").Value = 835.0000200000001
$ws.Range("M17").Value = -213.000008
$ws.Range("N17").Value = -1173.00002
$ws.Range("H26").Value = 32.5
$ws.Range("I26").Value = 32.5
$ws.Range("K26").Value = 97.5
$ws.Range("M26").Value = 190.5
$ws.Range("H34").Value = 734.0769
$ws.Range("J34").Value = 1444.3334
$ws.Range("L34").Value = 4333.0002
$ws.Range("N34").Value = -4501.0002
$ws.Range("H37").Value = 1000000000
$ws.Range("J37").Value = 1000000000
$ws.Range("L37").Value = 3000000000
$ws.Range("N37").Value = -3000000224
$ws.Range("H39").Value = 2192.6365
$ws.Range("J39").Value = 4574.75
$ws.Range("L39").Value = 13724.25
$ws.Range("N39").Value = -14312.25
$ws.Range("H131").Value = 92626.5
$ws.Range("J131").Value = 183761.27
$ws.Range("L131").Value = 551283.8099999999
$ws.Range("N131").Value = -561363.8099999999
$ws.Range("H137").Value = 3841.9473
$ws.Range("I137").Value = 2381.25
$ws.Range("J137").Value = 4231.467
$ws.Range("K137").Value = 7143.75
$ws.Range("L137").Value = 12694.401
$ws.Range("M137").Value = -2043.75
$ws.Range("N137").Value = -22894.401
$ws.Range("H138").Value = 3933
$ws.Range("I138").Value = 3800
$ws.Range("K138").Value = 11400
$ws.Range("M138").Value = -6260
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 303
$ws.Range("I4").Value = 303
$ws.Range("K4").Value = 303
$ws.Range("M4").Value = -191
$ws.Range("H122").Value = 2612.6667
$ws.Range("I122").Value = 2715
$ws.Range("J122").Value = 2232.5715
$ws.Range("K122").Value = 8145
$ws.Range("L122").Value = 6697.7145
$ws.Range("M122").Value = -5695
$ws.Range("N122").Value = -11597.7145
$ws.Range("H132").Value = 1594.2354
$ws.Range("I132").Value = 1600.125
$ws.Range("K132").Value = 4800.375
$ws.Range("M132").Value = -2270.375
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 6129.5
$ws.Range("I35").Value = 1173
$ws.Range("J35").Value = 20999
$ws.Range("K35").Value = 1173
$ws.Range("L35").Value = 20999
$ws.Range("M35").Value = -837
$ws.Range("N35").Value = -21671
$ws.Range("H136").Value = 3491.6155
$ws.Range("I136").Value = 2239.1
$ws.Range("K136").Value = 6717.299999999999
$ws.Range("M136").Value = -4167.299999999999
